# Commit "update file with jgit": cell E8 on the "Rules" sheet changes
# from "Good Morning" to "GIT UPDATE", and that cell becomes the active
# selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
